$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number that was updated
# from 45175 (2023-09-06) to 45177 (2023-09-08) for every data row
# (rows 2 through 471).
$ws.Range("C2:C471").Value = 45177
